# Apply cryptos list update (prices/volumes) matching source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "42.914.70"
Set-TextCell "E2" "  +0.01%  "

Set-TextCell "D3" "2.292.42"
Set-TextCell "E3" "  +1.71%  "

Set-TextCell "E4" "  -0.08%  "

Set-TextCell "D5" "252.05"
Set-TextCell "E5" "  +1.08%  "

Set-TextCell "D6" "0.637"
Set-TextCell "E6" "  +0.44%  "

Set-TextCell "D7" "75.35"
Set-TextCell "E7" "  +6.44%  "

Set-TextCell "E8" "  +0.00%  "

Set-TextCell "D9" "0.651"
Set-TextCell "E9" "  -3.60%  "

Set-TextCell "D10" "39.04"
Set-TextCell "E10" "  -0.40%  "

Set-TextCell "D11" "0.0985"
Set-TextCell "E11" "  +1.34%  "

Set-TextCell "D12" "7.56"
Set-TextCell "E12" "  -0.08%  "

Set-TextCell "D13" "0.107"
Set-TextCell "E13" "  +1.55%  "

Set-TextCell "D14" "2.635.58"
Set-TextCell "E14" "  +1.71%  "

Set-TextCell "D15" "15.14"
Set-TextCell "E15" "  +2.32%  "

Set-TextCell "D16" "0.873"
Set-TextCell "E16" "  -1.41%  "

Set-TextCell "D17" "2.290.01"
Set-TextCell "E17" "  +1.14%  "

Set-TextCell "D18" "42.791.21"
Set-TextCell "E18" "  +0.05%  "

Set-TextCell "E19" "  +1.91%  "

Set-TextCell "D20" "6.23"
Set-TextCell "E20" "  -0.57%  "

Set-TextCell "D21" "72.30"
Set-TextCell "E21" "  -1.11%  "

Set-TextCell "D22" "237.40"
Set-TextCell "E22" "  +0.72%  "

Set-TextCell "D23" "2.16"
Set-TextCell "E23" "  +4.98%  "

Set-TextCell "D24" "3.87"
Set-TextCell "E24" "  -2.15%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell "D25" "1.00"
Set-TextCell "E25" "  +0.11%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D26" "11.36"
Set-TextCell "E26" "  -1.76%  "

Set-TextCell "D27" "2.40"
Set-TextCell "E27" "  -1.13%  "

Set-TextCell "D28" "2.18"
Set-TextCell "E28" "  +1.46%  "

Set-TextCell "D29" "167.43"
Set-TextCell "E29" "  +0.13%  "

Set-TextCell "D30" "21.02"
Set-TextCell "E30" "  +0.02%  "

Set-TextCell "D31" "0.0861"
Set-TextCell "E31" "  +8.88%  "

Set-TextCell "D32" "6.22"
Set-TextCell "E32" "  -4.28%  "

Set-TextCell "D33" "0.127"
Set-TextCell "E33" "  -1.29%  "

Set-TextCell "D34" "31.49"
Set-TextCell "E34" "  +1.58%  "

Set-TextCell "E35" "  +1.14%  "

Set-TextCell "E36" "  +8.10%  "

Set-TextCell "D37" "4.79"
Set-TextCell "E37" "  +2.08%  "

Set-TextCell "D38" "0.0304"
Set-TextCell "E38" "  -5.46%  "

Set-TextCell "D39" "13.60"
Set-TextCell "E39" "  +8.58%  "

Set-TextCell "D40" "2.30"
Set-TextCell "E40" "  -0.01%  "

Set-TextCell "D41" "5.98"
Set-TextCell "E41" "  +2.87%  "

Set-TextCell "D42" "0.210"
Set-TextCell "E42" "  +4.04%  "

Set-TextCell "E43" "  +1.72%  "

Set-TextCell "D44" "61.08"
Set-TextCell "E44" "  -3.02%  "

Set-TextCell "D45" "4.85"
Set-TextCell "E45" "  -1.22%  "

Set-TextCell "D46" "105.79"
Set-TextCell "E46" "  +11.81%  "

Set-TextCell "D47" "0.101"
Set-TextCell "E47" "  -1.95%  "

Set-TextCell "E48" "  -0.06%  "

Set-TextCell "D49" "1.17"
Set-TextCell "E49" "  -0.65%  "

Set-TextCell "E50" "  -1.51%  "

Set-TextCell "D51" "4.23"
Set-TextCell "E51" "  -1.69%  "
